$wb = $excel.ActiveWorkbook

# ALC row 29
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 500.75
$ws.Range("I29").Value = 201.2
$ws.Range("J29").Value = 1000
$ws.Range("K29").Value = 603.5999999999999
$ws.Range("L29").Value = 3000
$ws.Range("M29").Value = -322.5999999999999
$ws.Range("N29").Value = -3562

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 310.33334
$ws.Range("I33").Value = 312.04166
$ws.Range("J33").Value = 296.66666
$ws.Range("K33").Value = 312.04166
$ws.Range("L33").Value = 296.66666
$ws.Range("M33").Value = -83.04165999999998
$ws.Range("N33").Value = -754.66666

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4321.5
$ws.Range("I74").Value = 3937.375
$ws.Range("K74").Value = 3937.375
$ws.Range("M74").Value = -3001.375

# ALC row 75
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 49314
$ws.Range("J75").Value = 49314
$ws.Range("L75").Value = 49314
$ws.Range("N75").Value = -51186

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 4321.5
$ws.Range("I77").Value = 3937.375
$ws.Range("K77").Value = 19686.875
$ws.Range("M77").Value = -15006.875

# ALC row 78
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H78").Value = 49314
$ws.Range("J78").Value = 49314
$ws.Range("L78").Value = 147942
$ws.Range("N78").Value = -157302

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2150.3
$ws.Range("I116").Value = 2186.1428
$ws.Range("J116").Value = 2066.6667
$ws.Range("K116").Value = 2186.1428
$ws.Range("L116").Value = 2066.6667
$ws.Range("M116").Value = 1255.8572
$ws.Range("N116").Value = -8950.6667

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1559.5294
$ws.Range("I2").Value = 1274.9166
$ws.Range("J2").Value = 2242.6
$ws.Range("K2").Value = 1274.9166
$ws.Range("L2").Value = 2242.6
$ws.Range("M2").Value = -1161.9166
$ws.Range("N2").Value = -2468.6

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11546.8955
$ws.Range("I32").Value = 11366.267
$ws.Range("K32").Value = 11366.267
$ws.Range("M32").Value = -11079.267

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 12822345
$ws.Range("I61").Value = 13159683
$ws.Range("J61").Value = 3500
$ws.Range("K61").Value = 13159683
$ws.Range("L61").Value = 3500
$ws.Range("M61").Value = -13159471
$ws.Range("N61").Value = -3924

# ARM row 86
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# ARM row 89
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1559.5294
$ws.Range("I116").Value = 1274.9166
$ws.Range("J116").Value = 2242.6
$ws.Range("K116").Value = 1274.9166
$ws.Range("L116").Value = 2242.6
$ws.Range("M116").Value = 1019.0834
$ws.Range("N116").Value = -6830.6

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 12822345
$ws.Range("I136").Value = 13159683
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 39479049
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -39476499
$ws.Range("N136").Value = -15600

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1559.5294
$ws.Range("I3").Value = 1274.9166
$ws.Range("J3").Value = 2242.6
$ws.Range("K3").Value = 1274.9166
$ws.Range("L3").Value = 2242.6
$ws.Range("M3").Value = -1160.9166
$ws.Range("N3").Value = -2470.6

# BSM row 92
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 25000
$ws.Range("J92").Value = 25000
$ws.Range("L92").Value = 25000
$ws.Range("N92").Value = -29992

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1181.1177
$ws.Range("I94").Value = 1119.9231
$ws.Range("J94").Value = 1380
$ws.Range("K94").Value = 1119.9231
$ws.Range("L94").Value = 1380
$ws.Range("M94").Value = -668.9231
$ws.Range("N94").Value = -2282

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 880.6667
$ws.Range("I99").Value = 846.36365
$ws.Range("J99").Value = 975
$ws.Range("K99").Value = 846.36365
$ws.Range("L99").Value = 975
$ws.Range("M99").Value = 651.63635
$ws.Range("N99").Value = -3971

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6249.0527
$ws.Range("I134").Value = 6376
$ws.Range("J134").Value = 6215.2
$ws.Range("K134").Value = 19128
$ws.Range("L134").Value = 18645.6
$ws.Range("M134").Value = -16593
$ws.Range("N134").Value = -23715.6

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7095926.5
$ws.Range("J31").Value = 37041156
$ws.Range("L31").Value = 37041156
$ws.Range("N31").Value = -37041746

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7095926.5
$ws.Range("J34").Value = 37041156
$ws.Range("L34").Value = 37041156
$ws.Range("N34").Value = -37041560

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 675.2273
$ws.Range("I107").Value = 436.45456
$ws.Range("K107").Value = 436.45456
$ws.Range("M107").Value = 1483.54544

# CUL row 80
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2196.1538
$ws.Range("I80").Value = 966.3333
$ws.Range("J80").Value = 2565.1
$ws.Range("K80").Value = 2898.9999
$ws.Range("L80").Value = 7695.299999999999
$ws.Range("M80").Value = -1962.9999
$ws.Range("N80").Value = -9567.299999999999

# CUL row 83
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 2196.1538
$ws.Range("I83").Value = 966.3333
$ws.Range("J83").Value = 2565.1
$ws.Range("K83").Value = 8696.9997
$ws.Range("L83").Value = 23085.9
$ws.Range("M83").Value = -4016.9997
$ws.Range("N83").Value = -32445.9

# CUL row 112
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 4748.95
$ws.Range("I112").Value = 4449.6665
$ws.Range("J112").Value = 4801.7646
$ws.Range("K112").Value = 13348.9995
$ws.Range("L112").Value = 14405.2938
$ws.Range("M112").Value = -12240.9995
$ws.Range("N112").Value = -16621.2938

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 860.45
$ws.Range("I113").Value = 622.8570999999999
$ws.Range("J113").Value = 988.38464
$ws.Range("K113").Value = 1868.5713
$ws.Range("L113").Value = 2965.15392
$ws.Range("M113").Value = 301.4287000000002
$ws.Range("N113").Value = -7305.15392

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 4274.294
$ws.Range("I129").Value = 1375.7142
$ws.Range("J129").Value = 6303.3
$ws.Range("K129").Value = 4127.142599999999
$ws.Range("L129").Value = 18909.9
$ws.Range("M129").Value = 872.8574000000008
$ws.Range("N129").Value = -28909.9

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 4894
$ws.Range("I140").Value = 2823.3333
$ws.Range("J140").Value = 8000
$ws.Range("K140").Value = 8469.999899999999
$ws.Range("L140").Value = 24000
$ws.Range("M140").Value = -3289.999899999999
$ws.Range("N140").Value = -34360

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1307.3182
$ws.Range("I107").Value = 1889.3077
$ws.Range("J107").Value = 466.66666
$ws.Range("K107").Value = 1889.3077
$ws.Range("L107").Value = 466.66666
$ws.Range("M107").Value = 30.69229999999993
$ws.Range("N107").Value = -4306.66666

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3511299.8
$ws.Range("I122").Value = 11113580
$ws.Range("J122").Value = 2555
$ws.Range("K122").Value = 33340740
$ws.Range("L122").Value = 7665
$ws.Range("M122").Value = -33338290
$ws.Range("N122").Value = -12565

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2603.2144
$ws.Range("I16").Value = 2620.4167
$ws.Range("J16").Value = 2500
$ws.Range("K16").Value = 2620.4167
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = -2450.4167
$ws.Range("N16").Value = -2840

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 171.16667
$ws.Range("I55").Value = 68.5
$ws.Range("J55").Value = 222.5
$ws.Range("K55").Value = 68.5
$ws.Range("L55").Value = 222.5
$ws.Range("M55").Value = 104.5
$ws.Range("N55").Value = -568.5

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 580.6
$ws.Range("I93").Value = 580.6
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 580.6
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 667.4
$ws.Range("N93").ClearContents()

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1000
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1000
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 721.5143
$ws.Range("I136").Value = 723.90625
$ws.Range("J136").Value = 696
$ws.Range("K136").Value = 2171.71875
$ws.Range("L136").Value = 2088
$ws.Range("M136").Value = 378.28125
$ws.Range("N136").Value = -7188
